$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new weekly entry for row 15
$ws.Range("B15").Value = 7.25
$ws.Range("C15").Value = "preprocessing redo"

# Update the selection to reflect where the user left off
$ws.Range("B15").Select()
